$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2954.9167
$ws.Range("J19").Value = 2940.4443
$ws.Range("L19").Value = 2940.4443
$ws.Range("N19").Value = -3290.4443
$ws.Range("H116").Value = 4713.7856
$ws.Range("I116").Value = 4049.25
$ws.Range("J116").Value = 5599.8335
$ws.Range("K116").Value = 4049.25
$ws.Range("L116").Value = 5599.8335
$ws.Range("M116").Value = -607.25
$ws.Range("N116").Value = -12483.8335
$ws.Range("H132").Value = 1081.0975
$ws.Range("I132").Value = 993.2222
$ws.Range("J132").Value = 1713.8
$ws.Range("K132").Value = 2979.6666
$ws.Range("L132").Value = 5141.4
$ws.Range("M132").Value = -449.6666
$ws.Range("N132").Value = -10201.4
$ws.Range("H138").Value = 2307.9348
$ws.Range("J138").Value = 2414.543
$ws.Range("L138").Value = 7243.629000000001
$ws.Range("N138").Value = -17523.629

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2014.12
$ws.Range("I2").Value = 944.6667
$ws.Range("K2").Value = 944.6667
$ws.Range("M2").Value = -831.6667
$ws.Range("H110").Value = 3504.2666
$ws.Range("I110").Value = 2097.3333
$ws.Range("K110").Value = 2097.3333
$ws.Range("M110").Value = -52.33329999999978
$ws.Range("H116").Value = 2014.12
$ws.Range("I116").Value = 944.6667
$ws.Range("K116").Value = 944.6667
$ws.Range("M116").Value = 1349.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2014.12
$ws.Range("I3").Value = 944.6667
$ws.Range("K3").Value = 944.6667
$ws.Range("M3").Value = -830.6667
$ws.Range("H26").Value = 15952.429
$ws.Range("I26").Value = 15952.429
$ws.Range("K26").Value = 15952.429
$ws.Range("M26").Value = -15660.429
$ws.Range("H57").Value = 199999
$ws.Range("J57").Value = 199999
$ws.Range("L57").Value = 199999
$ws.Range("N57").Value = -201439
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null
$ws.Range("H134").Value = 2822.2449
$ws.Range("I134").Value = 2899.8086
$ws.Range("K134").Value = 8699.425799999999
$ws.Range("M134").Value = -6164.425799999999
$ws.Range("H136").Value = 199999
$ws.Range("J136").Value = 199999
$ws.Range("L136").Value = 199999
$ws.Range("N136").Value = -210199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 604
$ws.Range("J22").Value = 926.3333
$ws.Range("L22").Value = 926.3333
$ws.Range("N22").Value = -1626.3333
$ws.Range("H31").Value = 3724.7942
$ws.Range("I31").Value = 3128.111
$ws.Range("K31").Value = 3128.111
$ws.Range("M31").Value = -2833.111
$ws.Range("H34").Value = 3724.7942
$ws.Range("I34").Value = 3128.111
$ws.Range("K34").Value = 3128.111
$ws.Range("M34").Value = -2926.111
$ws.Range("H60").Value = 20495.75
$ws.Range("H86").Value = 58072.8
$ws.Range("J86").Value = 78858.28999999999
$ws.Range("L86").Value = 78858.28999999999
$ws.Range("N86").Value = -81104.28999999999
$ws.Range("H89").Value = 58072.8
$ws.Range("J89").Value = 78858.28999999999
$ws.Range("L89").Value = 394291.45
$ws.Range("N89").Value = -405523.45
$ws.Range("H94").Value = 1157.75
$ws.Range("I94").Value = 1059
$ws.Range("J94").Value = 1228.2858
$ws.Range("K94").Value = 1059
$ws.Range("L94").Value = 1228.2858
$ws.Range("M94").Value = -608
$ws.Range("N94").Value = -2130.2858
$ws.Range("H99").Value = 5403.4546
$ws.Range("I99").Value = 4969.7144
$ws.Range("K99").Value = 4969.7144
$ws.Range("M99").Value = -3471.7144
$ws.Range("H122").Value = 3568.75
$ws.Range("J122").Value = 3591.4546
$ws.Range("L122").Value = 10774.3638
$ws.Range("N122").Value = -15674.3638
$ws.Range("H126").Value = 5403.4546
$ws.Range("I126").Value = 4969.7144
$ws.Range("K126").Value = 14909.1432
$ws.Range("M126").Value = -12439.1432
$ws.Range("H132").Value = 2377.075
$ws.Range("I132").Value = 2137.1482
$ws.Range("K132").Value = 6411.444600000001
$ws.Range("M132").Value = -3881.444600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 300
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null
$ws.Range("H125").Value = 17662.666
$ws.Range("I125").Value = 12988
$ws.Range("K125").Value = 38964
$ws.Range("M125").Value = -34044
$ws.Range("H126").Value = 11499.833
$ws.Range("I126").Value = 7999.6665
$ws.Range("K126").Value = 23998.9995
$ws.Range("M126").Value = -19058.9995
$ws.Range("H129").Value = 24599842
$ws.Range("J129").Value = 2010500
$ws.Range("L129").Value = 6031500
$ws.Range("N129").Value = -6041500
$ws.Range("H131").Value = 29414158
$ws.Range("I131").Value = 125000820
$ws.Range("J131").Value = 2878.2307
$ws.Range("K131").Value = 375002460
$ws.Range("L131").Value = 8634.6921
$ws.Range("M131").Value = -374997420
$ws.Range("N131").Value = -18714.6921
$ws.Range("H141").Value = 2483
$ws.Range("I141").Value = 2483
$ws.Range("K141").Value = 7449
$ws.Range("M141").Value = -2269

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 50000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 50000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 50000
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = -50226
$ws.Range("H16").Value = 50000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 50000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 50000
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = -50500
$ws.Range("H44").Value = 19995
$ws.Range("J44").Value = 19995
$ws.Range("L44").Value = 19995
$ws.Range("N44").Value = -21187
$ws.Range("H132").Value = 1878.2778
$ws.Range("I132").Value = 1729.9706
$ws.Range("K132").Value = 5189.9118
$ws.Range("M132").Value = -2659.9118

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 2000000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 2000000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 2000000
$ws.Range("M19").Value = $null
$ws.Range("N19").Value = -2000340
$ws.Range("H46").Value = 11536.257
$ws.Range("I46").Value = 4535.857
$ws.Range("J46").Value = 16203.19
$ws.Range("K46").Value = 4535.857
$ws.Range("L46").Value = 16203.19
$ws.Range("M46").Value = -4347.857
$ws.Range("N46").Value = -16579.19
$ws.Range("H96").Value = 39000
$ws.Range("J96").Value = 39000
$ws.Range("L96").Value = 39000
$ws.Range("N96").Value = -44492
$ws.Range("H122").Value = 3779.818
$ws.Range("I122").Value = 3368.4285
$ws.Range("K122").Value = 10105.2855
$ws.Range("M122").Value = -7655.2855
$ws.Range("H136").Value = 3926.862
$ws.Range("I136").Value = 3545.25
$ws.Range("K136").Value = 10635.75
$ws.Range("M136").Value = -8085.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9496
$ws.Range("J15").Value = 9495
$ws.Range("L15").Value = 9495
$ws.Range("N15").Value = -10071
$ws.Range("H122").Value = 4835.5454
$ws.Range("I122").Value = 4621.1875
$ws.Range("K122").Value = 13863.5625
$ws.Range("M122").Value = -11413.5625
$ws.Range("H132").Value = 4257.8965
$ws.Range("I132").Value = 3907.5
$ws.Range("J132").Value = 5939.8
$ws.Range("K132").Value = 11722.5
$ws.Range("L132").Value = 17819.4
$ws.Range("M132").Value = -9192.5
$ws.Range("N132").Value = -22879.4
